$wb = $excel.ActiveWorkbook

# Weekly data refresh: add data for 2022-07-14 (updates 2022 year-to-date
# totals in column I, plus a couple of historical 2015 corrections in column B).
$sheetUpdates = @{
    'Citywide Totals' = @{ "I2" = 3665; "I3" = 3791; "B4" = 1656; "I4" = 892; "I5" = 349; "I6" = 4255; "B7" = 23288; "I7" = 12952 }
    'By Neighborhood' = @{ "I2" = 113; "I7" = 419; "I8" = 787; "I15" = 154; "I19" = 345; "I20" = 318; "I25" = 63; "I27" = 121; "I29" = 841; "I31" = 122; "I33" = 581; "I34" = 62; "I37" = 416; "I40" = 24; "I42" = 452; "I48" = 167; "I52" = 282; "I54" = 291; "B63" = 362; "I63" = 54; "I64" = 117; "I66" = 33; "I67" = 499; "I71" = 38; "I79" = 348; "I83" = 262; "I85" = 587; "I90" = 162; "I91" = 158; "I94" = 119; "I95" = 205; "I98" = 86; "I100" = 25; "B101" = 23288; "I101" = 12952 }
    'South Shore' = @{ "I2" = 153; "I7" = 587 }
    'Little Village' = @{ "I4" = 28; "I6" = 70; "I7" = 282 }
    'Austin' = @{ "I2" = 246; "I4" = 45; "I6" = 254; "I7" = 787 }
    'Auburn Gresham' = @{ "I2" = 145; "I3" = 124; "I7" = 419 }
    'Grand Crossing' = @{ "I3" = 127; "I4" = 29; "I7" = 416 }
    'North Lawndale' = @{ "I3" = 175; "I7" = 499 }
    'Gage Park' = @{ "I6" = 38; "I7" = 122 }
    'South Chicago' = @{ "I2" = 91; "I3" = 103; "I7" = 262 }
    'West Pullman' = @{ "I2" = 72; "I3" = 77; "I7" = 205 }
    'Garfield Park' = @{ "I3" = 212; "I4" = 28; "I7" = 581 }
    'Loop' = @{ "I6" = 149; "I7" = 291 }
    'Englewood' = @{ "I6" = 231; "I7" = 841 }
    'Chatham' = @{ "I3" = 99; "I6" = 94; "I7" = 345 }
    'Lake View' = @{ "I6" = 95; "I7" = 167 }
    'Humboldt Park' = @{ "I2" = 123; "I7" = 452 }
    'Washington Park' = @{ "I3" = 54; "I6" = 50; "I7" = 158 }
    'Roseland' = @{ "I2" = 102; "I3" = 112; "I4" = 20; "I6" = 103; "I7" = 348 }
    'Near South Side' = @{ "I4" = 6; "I6" = 43; "I7" = 117 }
    'Chicago Lawn' = @{ "I3" = 97; "I7" = 318 }
    'Wrigleyville' = @{ "I2" = 4; "I5" = 14; "I6" = 25 }
    'Garfield Ridge' = @{ "I2" = 22; "I3" = 22; "I7" = 62 }
    'West Loop' = @{ "I2" = 21; "I6" = 66; "I7" = 119 }
    'East Side' = @{ "I6" = 20; "I7" = 63 }
    'Brighton Park' = @{ "I3" = 36; "I6" = 55; "I7" = 154 }
    'Wicker Park' = @{ "I6" = 52; "I7" = 86 }
    'North Center' = @{ "I6" = 14; "I7" = 33 }
    'Albany Park' = @{ "I4" = 13; "I7" = 113 }
    'Edgewater' = @{ "I2" = 32; "I7" = 121 }
    'Washington Heights' = @{ "I5" = 6; "I6" = 56; "I7" = 162 }
    'Oakland' = @{ "I3" = 14; "I7" = 38 }
    'Hegewisch' = @{ "I3" = 11; "I7" = 24 }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
